$d = $word.ActiveDocument

# Locate the paragraph that contains the LOB1012 requirement text, then
# remove the two blank paragraphs and the copyright paragraph that follow
# it (through the end of the copyright paragraph), leaving the blank
# paragraph + page-break paragraph that originally trailed the copyright
# text directly after LOB1012.

$start = $null
$end = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*LOB1012*") {
        $start = $p.Range.End
    }
    elseif ($start -ne $null -and $t -like "*Contact: luizeleno@usp.br*") {
        $end = $p.Range.End
        break
    }
}

if ($start -ne $null -and $end -ne $null) {
    $r = $d.Range($start, $end)
    $r.Delete()
}
